$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value2 = 1.02
$ws.Cells.Item(2, 3).Value2 = 1.092089594101554
$ws.Cells.Item(2, 4).Value2 = 1.0894326989828
$ws.Cells.Item(2, 5).Value2 = 1.103869260218573
$ws.Cells.Item(2, 6).Value2 = 1.106636349846574
$ws.Cells.Item(2, 9).Value2 = 1.066386861664593
$ws.Cells.Item(2, 10).Value2 = 1.096912042885264
$ws.Cells.Item(2, 11).Value2 = 1.092082831726952
$ws.Cells.Item(2, 12).Value2 = 1.106482848336906
$ws.Cells.Item(2, 13).Value2 = 1.109243050618817
$ws.Cells.Item(2, 14).Value2 = 1.098469783803731
$ws.Cells.Item(3, 2).Value2 = 1.02
$ws.Cells.Item(3, 3).Value2 = 1.09341973596796
$ws.Cells.Item(3, 4).Value2 = 1.090487402298105
$ws.Cells.Item(3, 5).Value2 = 1.105137986909428
$ws.Cells.Item(3, 6).Value2 = 1.107864071504002
$ws.Cells.Item(3, 9).Value2 = 1.066832309876386
$ws.Cells.Item(3, 10).Value2 = 1.09790610285978
$ws.Cells.Item(3, 11).Value2 = 1.092956903440685
$ws.Cells.Item(3, 12).Value2 = 1.107573085821558
$ws.Cells.Item(3, 13).Value2 = 1.110292876704602
$ws.Cells.Item(3, 14).Value2 = 1.099465255457431
$ws.Cells.Item(4, 2).Value2 = 1.02
$ws.Cells.Item(4, 3).Value2 = 1.094280119606616
$ws.Cells.Item(4, 4).Value2 = 1.091169495579245
$ws.Cells.Item(4, 5).Value2 = 1.105958914639338
$ws.Cells.Item(4, 6).Value2 = 1.108658422282863
$ws.Cells.Item(4, 9).Value2 = 1.067119146783773
$ws.Cells.Item(4, 10).Value2 = 1.098548460557073
$ws.Cells.Item(4, 11).Value2 = 1.093521507381645
$ws.Cells.Item(4, 12).Value2 = 1.108277948250959
$ws.Cells.Item(4, 13).Value2 = 1.110971537889316
$ws.Cells.Item(4, 14).Value2 = 1.100108525376334
$ws.Cells.Item(5, 2).Value2 = 1.02
$ws.Cells.Item(5, 3).Value2 = 1.094641754300481
$ws.Cells.Item(5, 4).Value2 = 1.091456160597114
$ws.Cells.Item(5, 5).Value2 = 1.106304029306156
$ws.Cells.Item(5, 6).Value2 = 1.108992353303799
$ws.Cells.Item(5, 9).Value2 = 1.067239399324075
$ws.Cells.Item(5, 10).Value2 = 1.09881830248133
$ws.Cells.Item(5, 11).Value2 = 1.093758634168894
$ws.Cells.Item(5, 12).Value2 = 1.108574132224351
$ws.Cells.Item(5, 13).Value2 = 1.111256694193149
$ws.Cells.Item(5, 14).Value2 = 1.100378750507075
$ws.Cells.Item(6, 2).Value2 = 1.02
$ws.Cells.Item(6, 3).Value2 = 1.094702470302718
$ws.Cells.Item(6, 4).Value2 = 1.091504287861353
$ws.Cells.Item(6, 5).Value2 = 1.106361975445337
$ws.Cells.Item(6, 6).Value2 = 1.109048421027386
$ws.Cells.Item(6, 9).Value2 = 1.067259570714494
$ws.Cells.Item(6, 10).Value2 = 1.098863598148924
$ws.Cells.Item(6, 11).Value2 = 1.093798435201387
$ws.Cells.Item(6, 12).Value2 = 1.108623854655178
$ws.Cells.Item(6, 13).Value2 = 1.111304564283213
$ws.Cells.Item(6, 14).Value2 = 1.100424110499712
$ws.Cells.Item(7, 2).Value2 = 1.02
$ws.Cells.Item(7, 3).Value2 = 1.09428495206079
$ws.Cells.Item(7, 4).Value2 = 1.091173326351024
$ws.Cells.Item(7, 5).Value2 = 1.105963526088795
$ws.Cells.Item(7, 6).Value2 = 1.108662884338139
$ws.Cells.Item(7, 9).Value2 = 1.067120754914101
$ws.Cells.Item(7, 10).Value2 = 1.098552067002112
$ws.Cells.Item(7, 11).Value2 = 1.093524676793722
$ws.Cells.Item(7, 12).Value2 = 1.108281906424014
$ws.Cells.Item(7, 13).Value2 = 1.110975348760043
$ws.Cells.Item(7, 14).Value2 = 1.100112136942939
$ws.Cells.Item(8, 2).Value2 = 1.02
$ws.Cells.Item(8, 3).Value2 = 1.092539187123266
$ws.Cells.Item(8, 4).Value2 = 1.089789218571944
$ws.Cells.Item(8, 5).Value2 = 1.104298038501546
$ws.Cells.Item(8, 6).Value2 = 1.107051279282909
$ws.Cells.Item(8, 9).Value2 = 1.066537692988744
$ws.Cells.Item(8, 10).Value2 = 1.097248170517343
$ws.Cells.Item(8, 11).Value2 = 1.092378432033182
$ws.Cells.Item(8, 12).Value2 = 1.106851423447334
$ws.Cells.Item(8, 13).Value2 = 1.109597979272928
$ws.Cells.Item(8, 14).Value2 = 1.098806388775602
$ws.Cells.Item(9, 2).Value2 = 1.02
$ws.Cells.Item(9, 3).Value2 = 1.089460442927776
$ws.Cells.Item(9, 4).Value2 = 1.087347333537505
$ws.Cells.Item(9, 5).Value2 = 1.101362958525146
$ws.Cells.Item(9, 6).Value2 = 1.104210818314531
$ws.Cells.Item(9, 9).Value2 = 1.065499514302322
$ws.Cells.Item(9, 10).Value2 = 1.094943820595392
$ws.Cells.Item(9, 11).Value2 = 1.09035103123919
$ws.Cells.Item(9, 12).Value2 = 1.10432608762668
$ws.Cells.Item(9, 13).Value2 = 1.107165840912621
$ws.Cells.Item(9, 14).Value2 = 1.096498766412449
$ws.Cells.Item(10, 2).Value2 = 1.02
$ws.Cells.Item(10, 3).Value2 = 1.087406091570049
$ws.Cells.Item(10, 4).Value2 = 1.085717347245777
$ws.Cells.Item(10, 5).Value2 = 1.099405909760585
$ws.Cells.Item(10, 6).Value2 = 1.102316638502384
$ws.Cells.Item(10, 9).Value2 = 1.064800106713958
$ws.Cells.Item(10, 10).Value2 = 1.093402945084932
$ws.Cells.Item(10, 11).Value2 = 1.088994235198173
$ws.Cells.Item(10, 12).Value2 = 1.102639278953636
$ws.Cells.Item(10, 13).Value2 = 1.105540906171477
$ws.Cells.Item(10, 14).Value2 = 1.094955702682024
$ws.Cells.Item(11, 2).Value2 = 1.02
$ws.Cells.Item(11, 3).Value2 = 1.086516052814854
$ws.Cells.Item(11, 4).Value2 = 1.085011031516255
$ws.Cells.Item(11, 5).Value2 = 1.098558375769961
$ws.Cells.Item(11, 6).Value2 = 1.101496279935233
$ws.Cells.Item(11, 9).Value2 = 1.064495512726316
$ws.Cells.Item(11, 10).Value2 = 1.092734599326018
$ws.Cells.Item(11, 11).Value2 = 1.088405471216606
$ws.Cells.Item(11, 12).Value2 = 1.101908071285194
$ws.Cells.Item(11, 13).Value2 = 1.104836431264255
$ws.Cells.Item(11, 14).Value2 = 1.094286407795471
$ws.Cells.Item(12, 2).Value2 = 1.02
$ws.Cells.Item(12, 3).Value2 = 1.086185375317449
$ws.Cells.Item(12, 4).Value2 = 1.084748593951356
$ws.Cells.Item(12, 5).Value2 = 1.098243543161337
$ws.Cells.Item(12, 6).Value2 = 1.101191534593718
$ws.Cells.Item(12, 9).Value2 = 1.06438210938171
$ws.Cells.Item(12, 10).Value2 = 1.092486172655384
$ws.Cells.Item(12, 11).Value2 = 1.088186586249891
$ws.Cells.Item(12, 12).Value2 = 1.101636344099622
$ws.Cells.Item(12, 13).Value2 = 1.104574625097515
$ws.Cells.Item(12, 14).Value2 = 1.09403762833047
$ws.Cells.Item(13, 2).Value2 = 1.02
$ws.Cells.Item(13, 3).Value2 = 1.086256310328542
$ws.Cells.Item(13, 4).Value2 = 1.084804891375932
$ws.Cells.Item(13, 5).Value2 = 1.098311076809781
$ws.Cells.Item(13, 6).Value2 = 1.101256904797879
$ws.Cells.Item(13, 9).Value2 = 1.064406446731952
$ws.Cells.Item(13, 10).Value2 = 1.092539468902688
$ws.Cells.Item(13, 11).Value2 = 1.08823354653275
$ws.Cells.Item(13, 12).Value2 = 1.101694636161781
$ws.Cells.Item(13, 13).Value2 = 1.104630789459532
$ws.Cells.Item(13, 14).Value2 = 1.094091000264559
$ws.Cells.Item(14, 2).Value2 = 1.02
$ws.Cells.Item(14, 3).Value2 = 1.086488720511649
$ws.Cells.Item(14, 4).Value2 = 1.084989339999231
$ws.Cells.Item(14, 5).Value2 = 1.098532352048238
$ws.Cells.Item(14, 6).Value2 = 1.101471090169638
$ws.Cells.Item(14, 9).Value2 = 1.064486144154724
$ws.Cells.Item(14, 10).Value2 = 1.092714067849934
$ws.Cells.Item(14, 11).Value2 = 1.088387382042323
$ws.Cells.Item(14, 12).Value2 = 1.101885612778239
$ws.Cells.Item(14, 13).Value2 = 1.104814793008782
$ws.Cells.Item(14, 14).Value2 = 1.094265847162336
$ws.Cells.Item(15, 2).Value2 = 1.02
$ws.Cells.Item(15, 3).Value2 = 1.086631905719014
$ws.Cells.Item(15, 4).Value2 = 1.085102974156891
$ws.Cells.Item(15, 5).Value2 = 1.098668684190163
$ws.Cells.Item(15, 6).Value2 = 1.101603053103107
$ws.Cells.Item(15, 9).Value2 = 1.064535213402711
$ws.Cells.Item(15, 10).Value2 = 1.092821620983728
$ws.Cells.Item(15, 11).Value2 = 1.088482139695774
$ws.Cells.Item(15, 12).Value2 = 1.102003263250726
$ws.Cells.Item(15, 13).Value2 = 1.104928145993691
$ws.Cells.Item(15, 14).Value2 = 1.094373553033917
$ws.Cells.Item(16, 2).Value2 = 1.02
$ws.Cells.Item(16, 3).Value2 = 1.087465149849415
$ws.Cells.Item(16, 4).Value2 = 1.085764211932558
$ws.Cells.Item(16, 5).Value2 = 1.099462155048659
$ws.Cells.Item(16, 6).Value2 = 1.102371079268992
$ws.Cells.Item(16, 9).Value2 = 1.064820284721106
$ws.Cells.Item(16, 10).Value2 = 1.093447276829501
$ws.Cells.Item(16, 11).Value2 = 1.089033282744559
$ws.Cells.Item(16, 12).Value2 = 1.102687789540714
$ws.Cells.Item(16, 13).Value2 = 1.105587641362063
$ws.Cells.Item(16, 14).Value2 = 1.095000097382754
$ws.Cells.Item(17, 2).Value2 = 1.02
$ws.Cells.Item(17, 3).Value2 = 1.08798768779886
$ws.Cells.Item(17, 4).Value2 = 1.086178847910918
$ws.Cells.Item(17, 5).Value2 = 1.099959844477965
$ws.Cells.Item(17, 6).Value2 = 1.102852795408749
$ws.Cells.Item(17, 9).Value2 = 1.064998634104153
$ws.Cells.Item(17, 10).Value2 = 1.093839428140078
$ws.Cells.Item(17, 11).Value2 = 1.089378660732298
$ws.Cells.Item(17, 12).Value2 = 1.103116956501675
$ws.Cells.Item(17, 13).Value2 = 1.106001091383975
$ws.Cells.Item(17, 14).Value2 = 1.095392805593173
$ws.Cells.Item(18, 2).Value2 = 1.02
$ws.Cells.Item(18, 3).Value2 = 1.088292428215283
$ws.Cells.Item(18, 4).Value2 = 1.086420647951529
$ws.Cells.Item(18, 5).Value2 = 1.100250127019694
$ws.Cells.Item(18, 6).Value2 = 1.103133756428462
$ws.Cells.Item(18, 9).Value2 = 1.065102493916001
$ws.Cells.Item(18, 10).Value2 = 1.094068053883354
$ws.Cells.Item(18, 11).Value2 = 1.08957999221841
$ws.Cells.Item(18, 12).Value2 = 1.103367204512432
$ws.Cells.Item(18, 13).Value2 = 1.106242166261267
$ws.Cells.Item(18, 14).Value2 = 1.095621756011228
$ws.Cells.Item(19, 2).Value2 = 1.02
$ws.Cells.Item(19, 3).Value2 = 1.088396328939506
$ws.Cells.Item(19, 4).Value2 = 1.086503087082899
$ws.Cells.Item(19, 5).Value2 = 1.100349104037583
$ws.Cells.Item(19, 6).Value2 = 1.103229554266518
$ws.Cells.Item(19, 9).Value2 = 1.065137878908421
$ws.Cells.Item(19, 10).Value2 = 1.09414599086033
$ws.Cells.Item(19, 11).Value2 = 1.089648620482985
$ws.Cells.Item(19, 12).Value2 = 1.103452519497688
$ws.Cells.Item(19, 13).Value2 = 1.1063243524721
$ws.Cells.Item(19, 14).Value2 = 1.095699803667651
$ws.Cells.Item(20, 2).Value2 = 1.02
$ws.Cells.Item(20, 3).Value2 = 1.087931629325024
$ws.Cells.Item(20, 4).Value2 = 1.086134366594594
$ws.Cells.Item(20, 5).Value2 = 1.099906448301294
$ws.Cells.Item(20, 6).Value2 = 1.102801113486314
$ws.Cells.Item(20, 9).Value2 = 1.064979516327914
$ws.Cells.Item(20, 10).Value2 = 1.093797365359459
$ws.Cells.Item(20, 11).Value2 = 1.08934161754231
$ws.Cells.Item(20, 12).Value2 = 1.103070919017606
$ws.Cells.Item(20, 13).Value2 = 1.105956740768059
$ws.Cells.Item(20, 14).Value2 = 1.095350683078581
$ws.Cells.Item(21, 2).Value2 = 1.02
$ws.Cells.Item(21, 3).Value2 = 1.086420283673829
$ws.Cells.Item(21, 4).Value2 = 1.08493502671552
$ws.Cells.Item(21, 5).Value2 = 1.098467192601541
$ws.Cells.Item(21, 6).Value2 = 1.10140801870068
$ws.Cells.Item(21, 9).Value2 = 1.06446268253636
$ws.Cells.Item(21, 10).Value2 = 1.092662657616133
$ws.Cells.Item(21, 11).Value2 = 1.08834208663123
$ws.Cells.Item(21, 12).Value2 = 1.101829378371765
$ws.Cells.Item(21, 13).Value2 = 1.10476061223625
$ws.Cells.Item(21, 14).Value2 = 1.094214363920107
$ws.Cells.Item(22, 2).Value2 = 1.02
$ws.Cells.Item(22, 3).Value2 = 1.085469590594577
$ws.Cells.Item(22, 4).Value2 = 1.084180487556215
$ws.Cells.Item(22, 5).Value2 = 1.097562153946839
$ws.Cells.Item(22, 6).Value2 = 1.100531963452384
$ws.Cells.Item(22, 9).Value2 = 1.064136203017912
$ws.Cells.Item(22, 10).Value2 = 1.09194821790543
$ws.Cells.Item(22, 11).Value2 = 1.087712531054345
$ws.Cells.Item(22, 12).Value2 = 1.101048053165789
$ws.Cells.Item(22, 13).Value2 = 1.104007789060033
$ws.Cells.Item(22, 14).Value2 = 1.093498909623067
$ws.Cells.Item(23, 2).Value2 = 1.02
$ws.Cells.Item(23, 3).Value2 = 1.085973614882724
$ws.Cells.Item(23, 4).Value2 = 1.084580528000757
$ws.Cells.Item(23, 5).Value2 = 1.098041944537435
$ws.Cells.Item(23, 6).Value2 = 1.100996393094812
$ws.Cells.Item(23, 9).Value2 = 1.064309421023483
$ws.Cells.Item(23, 10).Value2 = 1.092327051926735
$ws.Cells.Item(23, 11).Value2 = 1.088046376391461
$ws.Cells.Item(23, 12).Value2 = 1.101462317510605
$ws.Cells.Item(23, 13).Value2 = 1.104406948662481
$ws.Cells.Item(23, 14).Value2 = 1.093878281632135
$ws.Cells.Item(24, 2).Value2 = 1.02
$ws.Cells.Item(24, 3).Value2 = 1.087956959881865
$ws.Cells.Item(24, 4).Value2 = 1.086154465941613
$ws.Cells.Item(24, 5).Value2 = 1.099930575769962
$ws.Cells.Item(24, 6).Value2 = 1.102824466371425
$ws.Cells.Item(24, 9).Value2 = 1.064988155348294
$ws.Cells.Item(24, 10).Value2 = 1.093816372057694
$ws.Cells.Item(24, 11).Value2 = 1.089358356141225
$ws.Cells.Item(24, 12).Value2 = 1.10309172161541
$ws.Cells.Item(24, 13).Value2 = 1.105976781160404
$ws.Cells.Item(24, 14).Value2 = 1.095369716768507
$ws.Cells.Item(25, 2).Value2 = 1.02
$ws.Cells.Item(25, 3).Value2 = 1.090256685469825
$ws.Cells.Item(25, 4).Value2 = 1.087978975146308
$ws.Cells.Item(25, 5).Value2 = 1.102121794581025
$ws.Cells.Item(25, 6).Value2 = 1.104945230958514
$ws.Cells.Item(25, 9).Value2 = 1.065769188511708
$ws.Cells.Item(25, 10).Value2 = 1.095540359282055
$ws.Cells.Item(25, 11).Value2 = 1.090876070959051
$ws.Cells.Item(25, 12).Value2 = 1.104979510979176
$ws.Cells.Item(25, 13).Value2 = 1.107795216502089
$ws.Cells.Item(25, 14).Value2 = 1.097096152252471
